# Editado planilha com filiais
# 1. Rename the "base" sheet to "Sheet1"
# 2. Add three more branches (BMA, RSD, VRE) that belong to group "G4",
#    matching the existing rows already on the sheet (AGR / G4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("base")

$ws.Name = "Sheet1"

$ws.Range("A3").Value = "BMA"
$ws.Range("B3").Value = "G4"

$ws.Range("A4").Value = "RSD"
$ws.Range("B4").Value = "G4"

$ws.Range("A5").Value = "VRE"
$ws.Range("B5").Value = "G4"
